$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 172.52389046390363
$ws.Range("C2").Value = 223.70860614204946
$ws.Range("D2").Value = 174.11376059760292
$ws.Range("E2").Value = 221.86949841826072

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 171.82536395782364
$ws.Range("C3").Value = 226.60467210110875
$ws.Range("D3").Value = 177.64809922808044
$ws.Range("E3").Value = 217.48874091687628

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
